$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1798245614035088
$ws.Range("C2").Value = 0.5657894736842105
$ws.Range("J2").Value = 0.008771929824561403
$ws.Range("P2").Value = 0.1403508771929824
$ws.Range("S2").Value = 0.1052631578947368
$ws.Range("C3").Value = 0.04477611940298507
$ws.Range("J3").Value = 0.007462686567164179
$ws.Range("P3").Value = 0.746268656716418
$ws.Range("S3").Value = 0.2014925373134328
$ws.Range("P4").Value = 0.696969696969697
$ws.Range("S4").Value = 0.303030303030303
$ws.Range("P5").Value = 0.75
$ws.Range("S5").Value = 0.25
$ws.Range("B6").Value = 0.04615384615384616
$ws.Range("D6").Value = 0.02051282051282051
$ws.Range("E6").Value = 0.005128205128205128
$ws.Range("F6").Value = 0.03076923076923077
$ws.Range("J6").Value = 0.2051282051282051
$ws.Range("O6").Value = 0.01025641025641026
$ws.Range("Q6").Value = 0.2512820512820513
$ws.Range("R6").Value = 0.08205128205128205
$ws.Range("S6").Value = 0.3487179487179487
$ws.Range("B7").Value = 0.1020408163265306
$ws.Range("D7").Value = 0.01360544217687075
$ws.Range("E7").Value = 0.006802721088435374
$ws.Range("F7").Value = 0.06122448979591837
$ws.Range("J7").Value = 0.08843537414965986
$ws.Range("O7").Value = 0.006802721088435374
$ws.Range("Q7").Value = 0.2312925170068027
$ws.Range("R7").Value = 0.06802721088435375
$ws.Range("S7").Value = 0.4217687074829932
$ws.Range("B8").Value = 0.09392265193370165
$ws.Range("D8").Value = 0.005524861878453038
$ws.Range("E8").Value = 0.002762430939226519
$ws.Range("F8").Value = 0.06077348066298342
$ws.Range("J8").Value = 0.06353591160220995
$ws.Range("O8").Value = 0.005524861878453038
$ws.Range("Q8").Value = 0.2071823204419889
$ws.Range("R8").Value = 0.1132596685082873
$ws.Range("S8").Value = 0.4475138121546962
$ws.Range("B9").Value = 0.04371584699453552
$ws.Range("D9").Value = 0.04918032786885246
$ws.Range("F9").Value = 0.07103825136612021
$ws.Range("J9").Value = 0.0546448087431694
$ws.Range("O9").Value = 0.01092896174863388
$ws.Range("Q9").Value = 0.2349726775956284
$ws.Range("R9").Value = 0.1038251366120219
$ws.Range("S9").Value = 0.4316939890710382
$ws.Range("B10").Value = 0.1049773755656109
$ws.Range("D10").Value = 0.01447963800904977
$ws.Range("E10").Value = 0.0009049773755656109
$ws.Range("F10").Value = 0.07963800904977375
$ws.Range("J10").Value = 0.08144796380090498
$ws.Range("O10").Value = 0.01176470588235294
$ws.Range("Q10").Value = 0.2153846153846154
$ws.Range("R10").Value = 0.08416289592760182
$ws.Range("S10").Value = 0.4072398190045249
$ws.Range("G11").Value = 0.1371681415929203
$ws.Range("J11").Value = 0.06637168141592921
$ws.Range("K11").Value = 0.1991150442477876
$ws.Range("L11").Value = 0.5707964601769911
$ws.Range("S11").Value = 0.02654867256637168
$ws.Range("G12").Value = 0.7633587786259542
$ws.Range("J12").Value = 0.183206106870229
$ws.Range("K12").Value = 0.007633587786259542
$ws.Range("L12").Value = 0.03053435114503817
$ws.Range("S12").Value = 0.01526717557251908
$ws.Range("G13").Value = 0.6129032258064516
$ws.Range("J13").Value = 0.2903225806451613
$ws.Range("S13").Value = 0.09677419354838709
$ws.Range("F15").Value = 0.02873563218390805
$ws.Range("H15").Value = 0.1379310344827586
$ws.Range("I15").Value = 0.09770114942528736
$ws.Range("J15").Value = 0.3908045977011494
$ws.Range("K15").Value = 0.06896551724137931
$ws.Range("M15").Value = 0.01724137931034483
$ws.Range("N15").Value = 0.005747126436781609
$ws.Range("O15").Value = 0.04022988505747126
$ws.Range("S15").Value = 0.2126436781609195
$ws.Range("F16").Value = 0.03311258278145696
$ws.Range("H16").Value = 0.1390728476821192
$ws.Range("I16").Value = 0.0728476821192053
$ws.Range("J16").Value = 0.4437086092715232
$ws.Range("K16").Value = 0.09271523178807947
$ws.Range("M16").Value = 0.01986754966887417
$ws.Range("O16").Value = 0.04635761589403974
$ws.Range("S16").Value = 0.152317880794702
$ws.Range("F17").Value = 0.02277904328018223
$ws.Range("H17").Value = 0.1685649202733485
$ws.Range("I17").Value = 0.07517084282460136
$ws.Range("J17").Value = 0.4396355353075171
$ws.Range("K17").Value = 0.07517084282460136
$ws.Range("M17").Value = 0.01366742596810934
$ws.Range("O17").Value = 0.07289293849658314
$ws.Range("S17").Value = 0.132118451025057
$ws.Range("F18").Value = 0.02247191011235955
$ws.Range("H18").Value = 0.1235955056179775
$ws.Range("I18").Value = 0.1179775280898876
$ws.Range("J18").Value = 0.4831460674157304
$ws.Range("K18").Value = 0.08426966292134831
$ws.Range("M18").Value = 0.01123595505617977
$ws.Range("O18").Value = 0.06741573033707865
$ws.Range("S18").Value = 0.0898876404494382
$ws.Range("F19").Value = 0.009777777777777778
$ws.Range("H19").Value = 0.1982222222222222
$ws.Range("I19").Value = 0.08977777777777778
$ws.Range("J19").Value = 0.424
$ws.Range("K19").Value = 0.0951111111111111
$ws.Range("M19").Value = 0.01511111111111111
$ws.Range("N19").Value = 0.0008888888888888889
$ws.Range("O19").Value = 0.06311111111111112
$ws.Range("S19").Value = 0.104
